# Auto-generated Excel COM-interop script
# Applies updated market-price snapshot values (and removes now-stale
# profit/column cells) across the ALC/ARM/CRP/CUL/GSM/LTW/WVR sheets,
# matching the scheduled market-data refresh described in the commit.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1200
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()
$ws.Range("H20").Value = 1490
$ws.Range("I20").Value = 987.5
$ws.Range("J20").Value = 3500
$ws.Range("K20").Value = 987.5
$ws.Range("L20").Value = 3500
$ws.Range("M20").Value = -757.5
$ws.Range("N20").Value = -3960
$ws.Range("H32").Value = 792.8570999999999
$ws.Range("J32").Value = 783.3333
$ws.Range("L32").Value = 783.3333
$ws.Range("N32").Value = -1435.3333
$ws.Range("H35").Value = 1490
$ws.Range("I35").Value = 987.5
$ws.Range("J35").Value = 3500
$ws.Range("K35").Value = 987.5
$ws.Range("L35").Value = 3500
$ws.Range("M35").Value = -608.5
$ws.Range("N35").Value = -4258
$ws.Range("H38").Value = 54.125
$ws.Range("I38").Value = 54.125
$ws.Range("K38").Value = 162.375
$ws.Range("M38").Value = 209.625
$ws.Range("H61").Value = 312
$ws.Range("I61").Value = 312
$ws.Range("K61").Value = 936
$ws.Range("M61").Value = -764
$ws.Range("H86").Value = 14499
$ws.Range("I86").Value = 17398.8
$ws.Range("K86").Value = 17398.8
$ws.Range("M86").Value = -16275.8
$ws.Range("H88").Value = 1193.8889
$ws.Range("J88").Value = 1439.8
$ws.Range("L88").Value = 1439.8
$ws.Range("N88").Value = -2251.8
$ws.Range("H89").Value = 14499
$ws.Range("I89").Value = 17398.8
$ws.Range("K89").Value = 86994
$ws.Range("M89").Value = -81378
$ws.Range("H91").Value = 1193.8889
$ws.Range("J91").Value = 1439.8
$ws.Range("L91").Value = 1439.8
$ws.Range("N91").Value = -4247.8
$ws.Range("H110").Value = 40500
$ws.Range("J110").Value = 40500
$ws.Range("L110").Value = 40500
$ws.Range("N110").Value = -48680
$ws.Range("H116").Value = 6106.5
$ws.Range("I116").Value = 5723.3335
$ws.Range("J116").Value = 6336.4
$ws.Range("K116").Value = 5723.3335
$ws.Range("L116").Value = 6336.4
$ws.Range("M116").Value = -2281.3335
$ws.Range("N116").Value = -13220.4
$ws.Range("H125").Value = 8290.714
$ws.Range("I125").Value = 4333.3335
$ws.Range("J125").Value = 11258.75
$ws.Range("K125").Value = 39000.0015
$ws.Range("L125").Value = 101328.75
$ws.Range("M125").Value = -36540.0015
$ws.Range("N125").Value = -106248.75
$ws.Range("H131").Value = 1365.3334
$ws.Range("J131").Value = 3200
$ws.Range("L131").Value = 9600
$ws.Range("N131").Value = -19680

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6492.7427
$ws.Range("I32").Value = 5213.1177
$ws.Range("K32").Value = 5213.1177
$ws.Range("M32").Value = -4926.1177

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 2934
$ws.Range("I3").Value = 556.6667
$ws.Range("K3").Value = 556.6667
$ws.Range("M3").Value = -443.6667
$ws.Range("H4").Value = 5933.3335
$ws.Range("I4").Value = 200
$ws.Range("J4").Value = 8800
$ws.Range("K4").Value = 200
$ws.Range("L4").Value = 8800
$ws.Range("M4").Value = -88
$ws.Range("N4").Value = -9024
$ws.Range("H7").Value = 164.25
$ws.Range("I7").Value = 103
$ws.Range("J7").Value = 250
$ws.Range("K7").Value = 103
$ws.Range("L7").Value = 250
$ws.Range("M7").Value = 10
$ws.Range("N7").Value = -476
$ws.Range("H31").Value = 2277.9
$ws.Range("I31").Value = 1938.1666
$ws.Range("J31").Value = 2787.5
$ws.Range("K31").Value = 1938.1666
$ws.Range("L31").Value = 2787.5
$ws.Range("M31").Value = -1643.1666
$ws.Range("N31").Value = -3377.5
$ws.Range("H34").Value = 2277.9
$ws.Range("I34").Value = 1938.1666
$ws.Range("J34").Value = 2787.5
$ws.Range("K34").Value = 1938.1666
$ws.Range("L34").Value = 2787.5
$ws.Range("M34").Value = -1736.1666
$ws.Range("N34").Value = -3191.5
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()
$ws.Range("H38").Value = 38
$ws.Range("I38").Value = 38
$ws.Range("K38").Value = 38
$ws.Range("M38").Value = 339
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H42").Value = 999
$ws.Range("I42").Value = 999
$ws.Range("K42").Value = 999
$ws.Range("M42").Value = -406
$ws.Range("H46").Value = 38
$ws.Range("I46").Value = 38
$ws.Range("K46").Value = 38
$ws.Range("M46").Value = 173
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 8748.429
$ws.Range("I6").Value = 248
$ws.Range("K6").Value = 744
$ws.Range("M6").Value = -631
$ws.Range("H12").Value = 375
$ws.Range("J12").Value = 461.16666
$ws.Range("L12").Value = 1383.49998
$ws.Range("N12").Value = -1729.49998
$ws.Range("H37").Value = 96998.336
$ws.Range("J37").Value = 96998.336
$ws.Range("L37").Value = 290995.008
$ws.Range("N37").Value = -291219.008
$ws.Range("H40").Value = 257.85715
$ws.Range("I40").Value = 166.25
$ws.Range("J40").Value = 380
$ws.Range("K40").Value = 665
$ws.Range("L40").Value = 1520
$ws.Range("M40").Value = -596
$ws.Range("N40").Value = -1658
$ws.Range("H107").Value = 1164
$ws.Range("I107").Value = 1029.8
$ws.Range("J107").Value = 1238.5555
$ws.Range("K107").Value = 3089.4
$ws.Range("L107").Value = 3715.6665
$ws.Range("M107").Value = -1169.4
$ws.Range("N107").Value = -7555.666499999999
$ws.Range("H128").Value = 339949.5
$ws.Range("I128").Value = 339949.5
$ws.Range("K128").Value = 1019848.5
$ws.Range("M128").Value = -1014868.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2722.4666
$ws.Range("I126").Value = 3209.4
$ws.Range("J126").Value = 2479
$ws.Range("K126").Value = 9628.200000000001
$ws.Range("L126").Value = 7437
$ws.Range("M126").Value = -7158.200000000001
$ws.Range("N126").Value = -12377

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5150
$ws.Range("I7").Value = 4475
$ws.Range("K7").Value = 4475
$ws.Range("M7").Value = -4363
$ws.Range("H40").Value = 2600.8667
$ws.Range("I40").Value = 2043.1666
$ws.Range("J40").Value = 4831.6665
$ws.Range("K40").Value = 2043.1666
$ws.Range("L40").Value = 4831.6665
$ws.Range("M40").Value = -1907.1666
$ws.Range("N40").Value = -5103.6665
$ws.Range("H126").Value = 5150
$ws.Range("I126").Value = 4475
$ws.Range("K126").Value = 13425
$ws.Range("M126").Value = -10955
$ws.Range("H136").Value = 2432.6667
$ws.Range("I136").Value = 2310.7646
$ws.Range("K136").Value = 6932.293799999999
$ws.Range("M136").Value = -4382.293799999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()
$ws.Range("H11").Value = 1000000000
$ws.Range("J11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("N11").ClearContents()

Write-Output "Updated market snapshot values across 8 sheets."